# Roll the GSC export window forward by one day:
#   - drop the oldest date (2025-11-10) from the "Chart" sheet
#   - shift all remaining rows up by one
#   - re-append the row that gets pushed off the bottom of the used range
#   - append a brand-new row for the newest crawl day (2026-02-09)
# (the "Table" sheet just references the trailing shared strings and needs
#  no direct edits - its string indices float automatically.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 91

# Remember the metric values of the final existing data row; its own row
# position is about to be vacated by the delete below, so its content would
# otherwise be lost.
$carryNonHttps = $ws.Cells.Item($lastRow, 2).Value()
$carryHttps = $ws.Cells.Item($lastRow, 3).Value()

# Drop the oldest date row (row 2); every row below shifts up by one,
# which is exactly the "yesterday's neighbor becomes today's value" shift
# seen across the whole Pages column.
$ws.Rows.Item(2).Delete()

# The row that used to be last (now blank, since nothing shifted into it)
# becomes the next day after the new last populated row.
$ws.Cells.Item($lastRow, 1).NumberFormat = "@"
$ws.Cells.Item($lastRow, 1).Value = "2026-02-08"
$ws.Cells.Item($lastRow, 1).ClearFormats()
$ws.Cells.Item($lastRow, 2).Value = $carryNonHttps
$ws.Cells.Item($lastRow, 3).Value = $carryHttps

# Append the newest day's row.
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026-02-09"
$ws.Cells.Item($newRow, 1).ClearFormats()
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 28
